# Apply the "rerunning SCTv2 corrected annotation after rechunking Kriegstein
# ref data" update to the "logs" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

# --- Mark the four integrated SingleR annotation runs (rows 19-22) as removed ---
$ws.Range("E19").Value = "removed"
$ws.Range("E20").Value = "removed"
$ws.Range("E21").Value = "removed"
$ws.Range("E22").Value = "removed"

# --- Row 31: the Kriegstein rechunking entry now points at the new chunk folder
#     and has been completed (TODO after: annotation, pseudotime) ---
$ws.Range("B31").Value = "chunks_25"
$ws.Range("G31").Value = "annotation, pseudotime"

# --- New rows 32-35: rerun SCTv2 corrected annotation (new/old, pre/post
#     selection) against the rechunked Kriegstein reference data ---
$ws.Range("A32").Value = "Rdata "
$ws.Range("B32").Value = "SingleR_RData_2022-06-15 09-10-25"
$ws.Range("C32").Value = "Kriegstein to SingleR"
$ws.Range("D32").Value = "SCTv2 corrected new post selection"
$ws.Range("F32").Value = "rerun SCTv2 corrected pipeline"
$ws.Range("G32").Value = "pseudotime"

$ws.Range("A33").Value = "Rdata "
$ws.Range("B33").Value = "SingleR_RData_2022-06-15 09-11-29"
$ws.Range("C33").Value = "Kriegstein to SingleR"
$ws.Range("D33").Value = "SCTv2 corrected old post selection"
$ws.Range("F33").Value = "rerun SCTv2 corrected pipeline"
$ws.Range("G33").Value = "pseudotime"

$ws.Range("A34").Value = "Rdata "
$ws.Range("B34").Value = "SingleR_RData_2022-06-15 09-13-25"
$ws.Range("C34").Value = "Kriegstein to SingleR"
$ws.Range("D34").Value = "SCTv2 corrected old selection"
$ws.Range("F34").Value = "rerun SCTv2 corrected pipeline"
$ws.Range("G34").Value = "pseudotime"

$ws.Range("A35").Value = "Rdata "
$ws.Range("B35").Value = "SingleR_RData_2022-06-15 09-14-33"
$ws.Range("C35").Value = "Kriegstein to SingleR"
$ws.Range("D35").Value = "SCTv2 corrected new selection"
$ws.Range("F35").Value = "rerun SCTv2 corrected pipeline"
$ws.Range("G35").Value = "pseudotime"

# Leave the selection where the author left off after entering the new rows
$ws.Range("G32").Select()
